$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 8561.7099999999991
$ws.Range("B16").Value = 8979.24
$ws.Range("C16").Value = 17.2
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = $true
$ws.Range("F16").Value = 4.6500000000000004
$ws.Range("G16").Value = 42626.545590277776
$ws.Range("G16").NumberFormat = "m/d/yy h:mm"
$ws.Range("H16").Value = $false
